try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # Expanded city aliases: From/To route table (columns A and C) plus the
    # related revenue / margin / passenger figures in columns E,G,I,K,M,O.
    $routes = @(
        @{ Row = 3;  From = "Atlanta";  To = "New York";      E = 3602000;  G = 0.0809; I = 955000;   K = 0.09;  M = 245; O = 65 }
        @{ Row = 4;  From = "New York"; To = "Washington";    E = 4674000;  G = 0.105;  I = 336000;   K = 0.03;  M = 222; O = 16 }
        @{ Row = 5;  From = "Chicago";  To = "New York";      E = 4674000;  G = 0.0804; I = 1536000;  K = 0.14;  M = 550; O = 43 }
        @{ Row = 6;  From = "New York"; To = "Philadelphia";  E = 12180000; G = 0.1427; I = -716000;  K = -0.07; M = 321; O = -25 }
        @{ Row = 7;  From = "New York"; To = "San Francisco"; E = 3221000;  G = 0.0629; I = 1088000;  K = 0.04;  M = 436; O = 21 }
        @{ Row = 8;  From = "New York"; To = "Phoneix";       E = 2782000;  G = 0.0723; I = 467000;   K = 0.1;   M = 674; O = 33 }
    )

    foreach ($route in $routes) {
        $r = $route.Row
        $ws.Range("A$r").Value = $route.From
        $ws.Range("C$r").Value = $route.To
        $ws.Range("E$r").Value = $route.E
        $ws.Range("G$r").Value = $route.G
        $ws.Range("I$r").Value = $route.I
        $ws.Range("K$r").Value = $route.K
        $ws.Range("M$r").Value = $route.M
        $ws.Range("O$r").Value = $route.O
    }

    # Tweak the From/To column widths to fit the newly expanded city names.
    $ws.Columns.Item(1).ColumnWidth = 9.140625
    $ws.Columns.Item(3).ColumnWidth = 12.2299089431763
}
catch {
    Write-Host "Failed to update dashboard: $_"
    throw
}
